# Restore the "Integer max" value for rule R30 (row 10) from 18 to 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C10").Value = 1
